# Update the "Förändrad" date column (C) for rows 2-43 from serial 45789 to 45790
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45789) {
        $cell.Value2 = 45790
    }
}
